$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 220
$ws.Range("F3").Value = 11
$ws.Range("F4").Value = 403
$ws.Range("F5").Value = 203
$ws.Range("F6").Value = 808
$ws.Range("F7").Value = 108
$ws.Range("F8").Value = 10291
$ws.Range("F9").Value = 58
$ws.Range("F10").Value = 3544
$ws.Range("F11").Value = 213
$ws.Range("F12").Value = 2457
$ws.Range("F13").Value = 38
$ws.Range("F14").Value = 2829
$ws.Range("F16").Value = 511
$ws.Range("F17").Value = 2187
$ws.Range("F18").Value = 47
$ws.Range("F19").Value = 99
$ws.Range("F20").Value = 30
$ws.Range("F21").Value = 392
$ws.Range("F22").Value = 23
$ws.Range("F23").Value = 155
$ws.Range("F24").Value = 320
$ws.Range("F25").Value = 277
$ws.Range("F26").Value = 235
$ws.Range("F27").Value = 618
$ws.Range("F28").Value = 1326
$ws.Range("F29").Value = 14
$ws.Range("F30").Value = 1261
$ws.Range("F31").Value = 107
$ws.Range("F32").Value = 133
$ws.Range("F33").Value = 242
$ws.Range("F34").Value = 3832
$ws.Range("F35").Value = 3222
$ws.Range("F36").Value = 34
$ws.Range("F38").Value = 1047
$ws.Range("F39").Value = 404
$ws.Range("F40").Value = 7
$ws.Range("F41").Value = 1295
$ws.Range("F42").Value = 107
$ws.Range("F43").Value = 112
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 30
$ws.Range("F46").Value = 42
$ws.Range("F47").Value = 15

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 180
$ws.Range("F8").Value = 10
$ws.Range("F15").Value = 36
$ws.Range("F16").Value = 181

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 756
$ws.Range("F3").Value = 990
$ws.Range("F4").Value = 129
$ws.Range("F5").Value = 2056

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 756
$ws.Range("F3").Value = 990
$ws.Range("F4").Value = 129
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 403
$ws.Range("F8").Value = 203
$ws.Range("F9").Value = 808
$ws.Range("F10").Value = 108
$ws.Range("F11").Value = 10292
$ws.Range("F12").Value = 58
$ws.Range("F13").Value = 3545
$ws.Range("F14").Value = 213
$ws.Range("F15").Value = 2457
$ws.Range("F16").Value = 38
$ws.Range("F18").Value = 511
$ws.Range("F19").Value = 2187
$ws.Range("F20").Value = 47
$ws.Range("F21").Value = 99
$ws.Range("F22").Value = 30
$ws.Range("F23").Value = 392
$ws.Range("F24").Value = 155
$ws.Range("F25").Value = 320
$ws.Range("F26").Value = 235
$ws.Range("F27").Value = 1326
$ws.Range("F28").Value = 14
$ws.Range("F29").Value = 1261
$ws.Range("F30").Value = 107
$ws.Range("F31").Value = 133
$ws.Range("F33").Value = 10
$ws.Range("F36").Value = 3222
$ws.Range("F37").Value = 1047
$ws.Range("F41").Value = 7
$ws.Range("F43").Value = 36
$ws.Range("F44").Value = 1295
$ws.Range("F45").Value = 107
$ws.Range("F46").Value = 73
$ws.Range("F47").Value = 30
$ws.Range("F48").Value = 15
$ws.Range("F49").Value = 181
